$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 4.35
$ws.Range("J3").Value = 2.25
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 4.6
$ws.Range("O3").Value = 1.91
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = 1.27
$ws.Range("U3").Value = 1.85
$ws.Range("V3").Value = 1.75
$ws.Range("W3").Value = 6.4
$ws.Range("X3").Value = 7.7
$ws.Range("Z3").Value = 13.5
$ws.Range("AA3").Value = 14.5
$ws.Range("AB3").Value = 29
$ws.Range("AD3").Value = 6.8
$ws.Range("AE3").Value = 17
$ws.Range("AF3").Value = 90
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 11
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 14.5
$ws.Range("AK3").Value = 70
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 50

# Row 5
$ws.Range("G5").Value = 2.47
$ws.Range("H5").Value = 2.95
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.15
$ws.Range("K5").Value = 1.93
$ws.Range("L5").Value = 3.5
$ws.Range("M5").Value = 1.45
$ws.Range("N5").Value = 2.37
$ws.Range("O5").Value = 2.32
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 3.9
$ws.Range("R5").Value = 1.17
$ws.Range("S5").Value = 1.47
$ws.Range("T5").Value = 2.32
$ws.Range("U5").Value = 1.98
$ws.Range("V5").Value = 1.65
$ws.Range("W5").Value = 6.3
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 27
$ws.Range("AA5").Value = 25
$ws.Range("AB5").Value = 45
$ws.Range("AC5").Value = 6.8
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 100
$ws.Range("AH5").Value = 7.3
$ws.Range("AJ5").Value = 11
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 45

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 2.82
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 3.35
$ws.Range("R6").Value = 1.26
$ws.Range("T6").Value = 2.57
$ws.Range("X6").Value = 10.75
$ws.Range("Y6").Value = 9.25
$ws.Range("AA6").Value = 19.5
$ws.Range("AD6").Value = 6.4
$ws.Range("AH6").Value = 8.5
$ws.Range("AI6").Value = 14
$ws.Range("AL6").Value = 25
$ws.Range("AM6").Value = 37

# Row 7
$ws.Range("G7").Value = 1.93
$ws.Range("H7").Value = 3.55
$ws.Range("I7").Value = 3.45
$ws.Range("J7").Value = 2.52
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 3.9
$ws.Range("M7").Value = 1.27
$ws.Range("N7").Value = 3.1
$ws.Range("O7").Value = 1.8
$ws.Range("P7").Value = 1.8
$ws.Range("Q7").Value = 2.87
$ws.Range("R7").Value = 1.31
$ws.Range("T7").Value = 2.6
$ws.Range("V7").Value = 1.9
$ws.Range("W7").Value = 7.5
$ws.Range("X7").Value = 9.25
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 15.5
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 14.5
$ws.Range("AF7").Value = 65
$ws.Range("AH7").Value = 10.5
$ws.Range("AI7").Value = 18.5
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 45
$ws.Range("AL7").Value = 30
$ws.Range("AM7").Value = 37

# Row 9
$ws.Range("G9").Value = 2.52
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 2.55
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 2.07
$ws.Range("L9").Value = 3.15
$ws.Range("M9").Value = 1.26
$ws.Range("N9").Value = 3.15
$ws.Range("O9").Value = 1.78
$ws.Range("P9").Value = 1.82
$ws.Range("Q9").Value = 2.8
$ws.Range("R9").Value = 1.33
$ws.Range("S9").Value = 1.38
$ws.Range("T9").Value = 2.6
$ws.Range("U9").Value = 1.62
$ws.Range("V9").Value = 2.02
$ws.Range("W9").Value = 9.25
$ws.Range("X9").Value = 13.5
$ws.Range("Z9").Value = 28
$ws.Range("AA9").Value = 20
$ws.Range("AB9").Value = 27
$ws.Range("AC9").Value = 10.5
$ws.Range("AD9").Value = 6.4
$ws.Range("AE9").Value = 13
$ws.Range("AF9").Value = 55
$ws.Range("AG9").Value = 400
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 13.5
$ws.Range("AJ9").Value = 9.5
$ws.Range("AK9").Value = 28
$ws.Range("AL9").Value = 21
$ws.Range("AM9").Value = 28

# Row 11
$ws.Range("H11").Value = 2.75
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 3.2
$ws.Range("K11").Value = 1.8
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.62
$ws.Range("R11").Value = 1.13
$ws.Range("S11").Value = 1.67
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53
$ws.Range("W11").Value = 5.5
$ws.Range("AA11").Value = 26
$ws.Range("AC11").Value = 5
$ws.Range("AK11").Value = 41
$ws.Range("AL11").Value = 41
$ws.Range("AN11").Value = 1.14

# Row 12
$ws.Range("H12").Value = 3
$ws.Range("K12").Value = 1.95
$ws.Range("M12").Value = 1.44
$ws.Range("N12").Value = 2.63
$ws.Range("O12").Value = 2.5
$ws.Range("P12").Value = 1.5
$ws.Range("R12").Value = 1.17
$ws.Range("AH12").Value = 6.5
$ws.Range("AL12").Value = 21
$ws.Range("AN12").Value = 1.1
$ws.Range("AO12").Value = 7
$ws.Range("AP12").Value = 1.88
$ws.Range("AQ12").Value = 1.98

# Row 13
$ws.Range("G13").Value = 2.85
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 2.27
$ws.Range("J13").Value = 3.4
$ws.Range("K13").Value = 2.12
$ws.Range("L13").Value = 2.87
$ws.Range("M13").Value = 1.26
$ws.Range("N13").Value = 3.5
$ws.Range("O13").Value = 1.78
$ws.Range("P13").Value = 1.93
$ws.Range("Q13").Value = 2.82
$ws.Range("R13").Value = 1.38
$ws.Range("S13").Value = 1.38
$ws.Range("T13").Value = 2.82
$ws.Range("U13").Value = 1.62
$ws.Range("V13").Value = 2.18
$ws.Range("W13").Value = 10.25
$ws.Range("X13").Value = 16
$ws.Range("Y13").Value = 10.25
$ws.Range("Z13").Value = 35
$ws.Range("AA13").Value = 23
$ws.Range("AB13").Value = 28
$ws.Range("AC13").Value = 7.7
$ws.Range("AD13").Value = 6.5
$ws.Range("AE13").Value = 12
$ws.Range("AF13").Value = 50
$ws.Range("AH13").Value = 9
$ws.Range("AI13").Value = 12.5
$ws.Range("AJ13").Value = 9
$ws.Range("AK13").Value = 24
$ws.Range("AL13").Value = 17.5
$ws.Range("AM13").Value = 24
$ws.Range("AO13").Value = 7.7

# Row 14
$ws.Range("G14").Value = 2.77
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 2.2
$ws.Range("J14").Value = 3.3
$ws.Range("K14").Value = 2.27
$ws.Range("L14").Value = 2.72
$ws.Range("M14").Value = 1.2
$ws.Range("N14").Value = 4.05
$ws.Range("O14").Value = 1.62
$ws.Range("P14").Value = 2.18
$ws.Range("Q14").Value = 2.42
$ws.Range("S14").Value = 1.31
$ws.Range("T14").Value = 3.15
$ws.Range("U14").Value = 1.53
$ws.Range("V14").Value = 2.35
$ws.Range("W14").Value = 11.75
$ws.Range("X14").Value = 16.5
$ws.Range("Y14").Value = 10.25
$ws.Range("Z14").Value = 35
$ws.Range("AA14").Value = 21
$ws.Range("AB14").Value = 25
$ws.Range("AC14").Value = 8.5
$ws.Range("AD14").Value = 7.2
$ws.Range("AE14").Value = 11.75
$ws.Range("AF14").Value = 40
$ws.Range("AH14").Value = 10.5
$ws.Range("AI14").Value = 13
$ws.Range("AJ14").Value = 9
$ws.Range("AK14").Value = 23
$ws.Range("AL14").Value = 16
$ws.Range("AM14").Value = 21
$ws.Range("AO14").Value = 8.5

# Row 15
$ws.Range("G15").Value = 2.18
$ws.Range("H15").Value = 3.45
$ws.Range("I15").Value = 2.85
$ws.Range("J15").Value = 2.8
$ws.Range("K15").Value = 2.18
$ws.Range("L15").Value = 3.45
$ws.Range("S15").Value = 1.38
$ws.Range("T15").Value = 2.82
$ws.Range("W15").Value = 8.5
$ws.Range("X15").Value = 11
$ws.Range("Y15").Value = 9
$ws.Range("Z15").Value = 21
$ws.Range("AA15").Value = 17
$ws.Range("AB15").Value = 26
$ws.Range("AH15").Value = 9.75
$ws.Range("AI15").Value = 15
$ws.Range("AJ15").Value = 10.5
$ws.Range("AK15").Value = 35
$ws.Range("AL15").Value = 23
$ws.Range("AM15").Value = 30

# Row 16
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 2.82
$ws.Range("I16").Value = 3.1
$ws.Range("J16").Value = 3.1
$ws.Range("K16").Value = 1.88
$ws.Range("L16").Value = 3.9
$ws.Range("M16").Value = 1.5
$ws.Range("N16").Value = 2.4
$ws.Range("O16").Value = 2.5
$ws.Range("P16").Value = 1.47
$ws.Range("Q16").Value = 4.45
$ws.Range("R16").Value = 1.17
$ws.Range("S16").Value = 1.57
$ws.Range("T16").Value = 2.27
$ws.Range("U16").Value = 2.05
$ws.Range("V16").Value = 1.7
$ws.Range("W16").Value = 6.1
$ws.Range("Y16").Value = 9.75
$ws.Range("Z16").Value = 25
$ws.Range("AA16").Value = 23
$ws.Range("AB16").Value = 40
$ws.Range("AC16").Value = 5.4
$ws.Range("AD16").Value = 5.7
$ws.Range("AE16").Value = 17.5
$ws.Range("AF16").Value = 110
$ws.Range("AH16").Value = 7.2
$ws.Range("AI16").Value = 15
$ws.Range("AJ16").Value = 11.75
$ws.Range("AK16").Value = 45
$ws.Range("AL16").Value = 35
$ws.Range("AM16").Value = 55
$ws.Range("AN16").Value = 1.12
$ws.Range("AO16").Value = 5.4
